$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 192, shifting existing rows 192:239 down to 193:240
$ws.Rows(192).Insert()

# Populate the newly inserted row 192 with the new record's data.
$ws.Range("A192").Value = 3
$ws.Range("B192").Value = "Femacal de La Calera"
$ws.Range("C192").Value = "Coquimbo"
$ws.Range("D192").Value = 44511
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 100112040
$ws.Range("G192").Value = "Cilantro"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 220
$ws.Range("K192").Value = 2000
$ws.Range("L192").Value = 2500
$ws.Range("M192").Value = 2273
$ws.Range("N192").Value = "$/docena de atados (3 kilos)"
$ws.Range("O192").Value = "Provincia de Quillota"
$ws.Range("P192").Value = 758
$ws.Range("Q192").Value = 3
$ws.Range("R192").Value = "Hortaliza"
